$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.372.37'
$ws.Range('E2').Value = '  -1.80%  '
$ws.Range('D3').Value = '1.819.27'
$ws.Range('E4').Value = '  -0.93%  '
$ws.Range('D5').Value = "'330.56"
$ws.Range('E5').Value = '  -1.93%  '
$ws.Range('D6').Value = "'1.003"
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('D7').Value = "'0.4570"
$ws.Range('E7').Value = '  -2.44%  '
$ws.Range('D8').Value = "'0.3807"
$ws.Range('E8').Value = '  -4.01%  '
$ws.Range('D9').Value = "'46.06"
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = "'0.07864"
$ws.Range('E10').Value = '  -2.02%  '
$ws.Range('D11').Value = "'0.9601"
$ws.Range('E11').Value = '  -5.36%  '
$ws.Range('D12').Value = "'20.96"
$ws.Range('E12').Value = '  -4.88%  '
$ws.Range('D13').Value = '1.854.45'
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').Value = "'5.851"
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('D15').Value = "'7.072"
$ws.Range('E15').Value = '  -3.08%  '
$ws.Range('D16').Value = "'1.004"
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').Value = "'89.06"
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = "'0.06589"
$ws.Range('E19').Value = '  -2.82%  '
$ws.Range('D20').Value = "'17.11"
$ws.Range('E20').Value = '  -1.67%  '
$ws.Range('D21').Value = "'1.003"
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('D22').Value = '27.356.89'
$ws.Range('E22').Value = '  -1.88%  '
$ws.Range('E23').Value = '  -3.94%  '
$ws.Range('D24').Value = "'10.81"
$ws.Range('E24').Value = '  -1.96%  '
$ws.Range('D25').Value = "'2.265"
$ws.Range('E25').Value = '  -2.22%  '
$ws.Range('D26').Value = '2.061.20'
$ws.Range('E26').Value = '  -2.44%  '
$ws.Range('D27').Value = "'155.76"
$ws.Range('E27').Value = '  -2.22%  '
$ws.Range('D28').Value = "'19.31"
$ws.Range('E28').Value = '  -2.80%  '
$ws.Range('D29').Value = "'2.038"
$ws.Range('E29').Value = '  -5.82%  '
$ws.Range('D30').Value = "'5.235"
$ws.Range('E30').Value = '  -5.08%  '
$ws.Range('D31').Value = "'117.74"
$ws.Range('E31').Value = '  -3.61%  '
$ws.Range('D32').Value = "'0.09304"
$ws.Range('E32').Value = '  -1.91%  '
$ws.Range('D33').Value = "'0.9312"
$ws.Range('D34').Value = "'3.566"
$ws.Range('E34').Value = '  -2.16%  '
$ws.Range('D35').Value = "'5.215"
$ws.Range('E35').Value = '  -2.54%  '
$ws.Range('D36').Value = "'1.309"
$ws.Range('E36').Value = '  -3.77%  '
$ws.Range('D38').Value = "'0.02170"
$ws.Range('E38').Value = '  -3.56%  '
$ws.Range('D39').Value = "'8.069"
$ws.Range('E39').Value = '  -3.64%  '
$ws.Range('E40').Value = '  -0.74%  '
$ws.Range('D41').Value = "'1.138"
$ws.Range('E41').Value = '  -5.42%  '
$ws.Range('D42').Value = "'0.5748"
$ws.Range('E42').Value = '  -4.34%  '
$ws.Range('D43').Value = "'0.1815"
$ws.Range('E43').Value = '  -4.43%  '
$ws.Range('D44').Value = "'9.900"
$ws.Range('E44').Value = '  -5.08%  '
$ws.Range('D45').Value = "'1.281"
$ws.Range('E45').Value = '  +2.55%  '
$ws.Range('D46').Value = "'11.86"
$ws.Range('E46').Value = '  -3.91%  '
$ws.Range('D47').Value = "'0.5394"
$ws.Range('E47').Value = '  -4.97%  '
$ws.Range('D48').Value = "'1.867"
$ws.Range('E48').Value = '  -4.36%  '
$ws.Range('D49').Value = "'0.06574"
$ws.Range('E49').Value = '  -3.20%  '
$ws.Range('D50').Value = "'109.61"
$ws.Range('E50').Value = '  -2.88%  '
$ws.Range('E51').Value = '  -34.16%  '
